# Adds three more "load-case" column triplets (c11, c12, c13) to the
# right of the existing data (which ends at column AI / c10), mirroring
# the pattern already present in the sheet: a bold/centered/bordered
# header row plus 10 data rows with repeating numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the formatting of the last existing
#     header cell (AI1) onto the 9 new header cells, then set their text.
$ws.Range("AI1").Copy()
$ws.Range("AJ1:AR1").PasteSpecial(-4122)   # xlPasteFormats

$headers = @("Fz-c11", "Mx-c11", "My-c11", "Fz-c12", "Mx-c12", "My-c12", "Fz-c13", "Mx-c13", "My-c13")
$cols = @("AJ", "AK", "AL", "AM", "AN", "AO", "AP", "AQ", "AR")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# --- Data rows (2-11): same 9 numeric values repeat on every row.
$values = @(22, 22, 33, 24, 24, 36, 26, 26, 39)

for ($row = 2; $row -le 11; $row++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
